$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 2-4
# Column B holds plain numbers; column C holds text that happens to look numeric,
# so a leading apostrophe forces it to stay text (matching the original inlineStr cells).
# Resetting the style back to "Normal" afterwards drops the quote-prefix flag Excel
# would otherwise persist, keeping the cell format identical to the rest of the sheet.
$ws.Range("B2").Value = 1

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "'1"
$ws.Range("C3").Style = "Normal"

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "'3"
$ws.Range("C4").Style = "Normal"

# Append two new rows (28 and 29) after the existing data (which ended at row 27)
$ws.Range("A28").Value = "'"
$ws.Range("A28").Style = "Normal"
$ws.Range("B28").Value = "NIM"
$ws.Range("C28").Value = "Nama"

$ws.Range("A29").Value = "'"
$ws.Range("A29").Style = "Normal"
$ws.Range("B29").Value = 6
$ws.Range("C29").Value = "'6"
$ws.Range("C29").Style = "Normal"
